$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1621.091
$ws.Range("I19").Value = 861.75
$ws.Range("J19").Value = 2055
$ws.Range("K19").Value = 861.75
$ws.Range("L19").Value = 2055
$ws.Range("M19").Value = -686.75
$ws.Range("N19").Value = -2405
$ws.Range("H76").Value = 3349848.5
$ws.Range("I76").Value = 4687273
$ws.Range("J76").Value = 6287.5
$ws.Range("K76").Value = 4687273
$ws.Range("L76").Value = 6287.5
$ws.Range("M76").Value = -4686958
$ws.Range("N76").Value = -6917.5
$ws.Range("H79").Value = 3349848.5
$ws.Range("I79").Value = 4687273
$ws.Range("J79").Value = 6287.5
$ws.Range("K79").Value = 4687273
$ws.Range("L79").Value = 6287.5
$ws.Range("M79").Value = -4686181
$ws.Range("N79").Value = -8471.5
$ws.Range("H81").Value = 32500
$ws.Range("J81").Value = 32500
$ws.Range("L81").Value = 32500
$ws.Range("N81").Value = -34496
$ws.Range("H84").Value = 32500
$ws.Range("J84").Value = 32500
$ws.Range("L84").Value = 97500
$ws.Range("N84").Value = -107484
$ws.Range("H99").Value = 1651.1666
$ws.Range("J99").Value = 2933
$ws.Range("L99").Value = 8799
$ws.Range("N99").Value = -11795
$ws.Range("H117").Value = 49999
$ws.Range("J117").Value = 49999
$ws.Range("L117").Value = 49999
$ws.Range("N117").Value = -59177
$ws.Range("H132").Value = 925.32434
$ws.Range("I132").Value = 925.32434
$ws.Range("K132").Value = 2775.97302
$ws.Range("M132").Value = -245.9730199999999
$ws.Range("H137").Value = 2040.8
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("H138").Value = 3274.48
$ws.Range("I138").Value = 4035.1333
$ws.Range("J138").Value = 2133.5
$ws.Range("K138").Value = 12105.3999
$ws.Range("L138").Value = 6400.5
$ws.Range("M138").Value = -6965.3999
$ws.Range("N138").Value = -16680.5
$ws.Range("M137").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3001.4443
$ws.Range("I32").Value = 1913.8572
$ws.Range("K32").Value = 1913.8572
$ws.Range("M32").Value = -1626.8572
$ws.Range("H45").Value = 22501340
$ws.Range("J45").Value = 1824.5
$ws.Range("L45").Value = 1824.5
$ws.Range("N45").Value = -2578.5
$ws.Range("H61").Value = 3594.348
$ws.Range("I61").Value = 2776.7222
$ws.Range("K61").Value = 2776.7222
$ws.Range("M61").Value = -2564.7222
$ws.Range("H96").Value = 49648
$ws.Range("J96").Value = 49648
$ws.Range("L96").Value = 49648
$ws.Range("N96").Value = -55140
$ws.Range("H105").Value = 49850
$ws.Range("J105").Value = 49850
$ws.Range("L105").Value = 49850
$ws.Range("N105").Value = -56838
$ws.Range("H122").Value = 1846.9565
$ws.Range("I122").Value = 1889.4706
$ws.Range("K122").Value = 5668.4118
$ws.Range("M122").Value = -3218.4118
$ws.Range("H136").Value = 3594.348
$ws.Range("I136").Value = 2776.7222
$ws.Range("K136").Value = 8330.1666
$ws.Range("M136").Value = -5780.1666

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("H64").Value = 528.1
$ws.Range("I64").Value = 557
$ws.Range("K64").Value = 557
$ws.Range("M64").Value = -332
$ws.Range("H67").Value = 528.1
$ws.Range("I67").Value = 557
$ws.Range("K67").Value = 557
$ws.Range("M67").Value = 223
$ws.Range("H99").Value = 1517.3572
$ws.Range("I99").Value = 1494
$ws.Range("J99").Value = 1559.4
$ws.Range("K99").Value = 1494
$ws.Range("L99").Value = 1559.4
$ws.Range("M99").Value = 4
$ws.Range("N99").Value = -4555.4
$ws.Range("N51").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4029.2856
$ws.Range("I31").Value = 1168.3334
$ws.Range("K31").Value = 1168.3334
$ws.Range("M31").Value = -873.3334
$ws.Range("H34").Value = 4029.2856
$ws.Range("I34").Value = 1168.3334
$ws.Range("K34").Value = 1168.3334
$ws.Range("M34").Value = -966.3334
$ws.Range("H93").Value = 4899
$ws.Range("I93").Value = 4899
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 4899
$ws.Range("L93").Value = 0
$ws.Range("H134").Value = 3237.3635
$ws.Range("I134").Value = 2745.7778
$ws.Range("J134").Value = 5449.5
$ws.Range("K134").Value = 8237.3334
$ws.Range("L134").Value = 16348.5
$ws.Range("M134").Value = -5702.3334
$ws.Range("N134").Value = -21418.5
$ws.Range("N93").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 1678.5
$ws.Range("J111").Value = 3030
$ws.Range("L111").Value = 9090
$ws.Range("H115").Value = 2950
$ws.Range("J115").Value = 3900
$ws.Range("L115").Value = 11700
$ws.Range("N115").Value = -14050
$ws.Range("H122").Value = 1218.4445
$ws.Range("I122").Value = 960
$ws.Range("K122").Value = 8640
$ws.Range("M122").Value = -6190
$ws.Range("H129").Value = 39172.367
$ws.Range("I129").Value = 846
$ws.Range("J129").Value = 56861.46
$ws.Range("K129").Value = 2538
$ws.Range("L129").Value = 170584.38
$ws.Range("M129").Value = 2462
$ws.Range("N129").Value = -180584.38
$ws.Range("H131").Value = 8487451
$ws.Range("J131").Value = 14266.227
$ws.Range("L131").Value = 42798.681
$ws.Range("N131").Value = -52878.681
$ws.Range("H132").Value = 1945.3636
$ws.Range("I132").Value = 1514.2858
$ws.Range("J132").Value = 2699.75
$ws.Range("K132").Value = 13628.5722
$ws.Range("L132").Value = 24297.75
$ws.Range("M132").Value = -11098.5722
$ws.Range("N132").Value = -29357.75
$ws.Range("H138").Value = 1706.7142
$ws.Range("I138").Value = 1533
$ws.Range("J138").Value = 2749
$ws.Range("K138").Value = 4599
$ws.Range("L138").Value = 8247
$ws.Range("M138").Value = 541
$ws.Range("N138").Value = -18527
$ws.Range("N111").Value = -15224

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3984
$ws.Range("I132").Value = 3386.9167
$ws.Range("K132").Value = 10160.7501
$ws.Range("M132").Value = -7630.750100000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1886.3125
$ws.Range("I132").Value = 1129.5385
$ws.Range("J132").Value = 5165.6665
$ws.Range("K132").Value = 3388.6155
$ws.Range("L132").Value = 15496.9995
$ws.Range("M132").Value = -858.6155000000003
$ws.Range("N132").Value = -20556.9995
$ws.Range("H136").Value = 3366.7878
$ws.Range("I136").Value = 2195.818
$ws.Range("J136").Value = 5708.727
$ws.Range("K136").Value = 6587.454000000001
$ws.Range("L136").Value = 17126.181
$ws.Range("M136").Value = -4037.454000000001
$ws.Range("N136").Value = -22226.181

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 18000
$ws.Range("J68").Value = 18000
$ws.Range("L68").Value = 18000
$ws.Range("H71").Value = 18000
$ws.Range("J71").Value = 18000
$ws.Range("L71").Value = 54000
$ws.Range("H81").Value = 183.16667
$ws.Range("I81").Value = 183.16667
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 366.33334
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 694.66666
$ws.Range("H84").Value = 183.16667
$ws.Range("I84").Value = 183.16667
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 1831.6667
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 3472.3333
$ws.Range("H132").Value = 1840.4667
$ws.Range("I132").Value = 1134.0834
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 3402.2502
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -872.2501999999999
$ws.Range("N132").Value = -19058
$ws.Range("H136").Value = 1939.449
$ws.Range("I136").Value = 1632.2
$ws.Range("K136").Value = 4896.6
$ws.Range("M136").Value = -2346.6
$ws.Range("N68").Value = -19622
$ws.Range("N71").Value = -62112
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()
